# Atualização de bases das ligas, do dia: 03-05-2024 às 22:15
# Israel Premier League sheet: rotate the match-odds rows 108-110 (columns B:AB)
# up by one position (108<-109, 109<-110, 110<-108) and refresh several
# Asian-handicap odds on rows 215, 216, 217, 218, 220 and 221.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Israel Premier League")

# ---------------------------------------------------------------------------
# Rows 108-110: cyclic shift of columns B:AB (column A, the running index,
# stays put). Capture the current values first so the rotation is safe
# regardless of write order.
# ---------------------------------------------------------------------------
$row108 = $ws.Range("B108:AB108").Value
$row109 = $ws.Range("B109:AB109").Value
$row110 = $ws.Range("B110:AB110").Value

$ws.Range("B108:AB108").Value = $row109
$ws.Range("B109:AB109").Value = $row110
$ws.Range("B110:AB110").Value = $row108

# ---------------------------------------------------------------------------
# Rows 215-221: updated Asian-handicap / over-under odds.
# ---------------------------------------------------------------------------
$ws.Range("Q215").Value = 2.025
$ws.Range("R215").Value = 1.825

$ws.Range("Q216").Value = 2.025
$ws.Range("R216").Value = 1.825
$ws.Range("T216").Value = 1.85
$ws.Range("U216").Value = 2

$ws.Range("T217").Value = 1.85
$ws.Range("U217").Value = 2

$ws.Range("M218").Value = 3
$ws.Range("O218").Value = 2.5
$ws.Range("Q218").Value = 2.125
$ws.Range("R218").Value = 1.75

$ws.Range("M220").Value = 2
$ws.Range("N220").Value = 3.5
$ws.Range("O220").Value = 3.4
$ws.Range("Q220").Value = 2.05
$ws.Range("R220").Value = 1.8
$ws.Range("T220").Value = 2.025
$ws.Range("U220").Value = 1.825

$ws.Range("M221").Value = 8
$ws.Range("O221").Value = 1.4
$ws.Range("Q221").Value = 1.9
$ws.Range("R221").Value = 1.95
$ws.Range("T221").Value = 1.95
$ws.Range("U221").Value = 1.9
